$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 34 (shifts existing rows 34-45 down to 35-46),
# matching the newly agreed metier code LLS_LPF_0_0_0_SWP added 2024-03-11.
$ws.Rows(34).Insert()

# Start from a clean slate for the new row's formatting, then apply the
# yellow highlight fill used elsewhere in the sheet for newly-added codes.
$ws.Range("A34:F34").ClearFormats()
$ws.Range("A34:F34").Interior.Color = 65535

# Column E keeps the slightly smaller "Arial 9" look used by the rest of
# the Description column - grab that formatting from an existing cell.
$ws.Range("E2").Copy()
$ws.Range("E34").PasteSpecial(-4122)
$ws.Range("E34").Interior.Color = 65535

# New row values. Order matters for how new entries land in the shared
# string table (shortest/most-generic code first), so write D, C, B, then
# the rest, to mirror the table layout in the target workbook.
$ws.Range("D34").Value = "LLS_LPF"
$ws.Range("C34").Value = "LLS_LPF_0_0_0"
$ws.Range("B34").Value = "LLS_LPF_0_0_0_SWP"
$ws.Range("A34").Value = "MBS"
$ws.Range("E34").Value = "Swordfish (Xiphias gladius)"
$ws.Range("F34").Value = 20240311

# Refresh the AutoFilter so it covers the new last row of the table.
$ws.AutoFilterMode = $false
[void]$ws.Range("A1:F46").AutoFilter()

# Keep the hidden _FilterDatabase defined name in sync with the new range.
$names = $wb.Names
$n = $names.Item(1)
$n.RefersTo = "=Sheet1!`$A`$1:`$F`$46"

# Match the saved selection position recorded in the workbook.
[void]$ws.Range("E30").Select()
